$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: "Resolution" -> "Resolution/Bin size"
$ws.Range("B1").Value = "Resolution/Bin size"

# Normalize tool-name capitalization: ARROWHEAD -> Arrowhead (rows 2-6)
$ws.Range("A2:A6").Value = "Arrowhead"

# Normalize tool-name capitalization: PEAKACHU -> Peakachu (row 7)
$ws.Range("A7").Value = "Peakachu"

# Update the active selection to match the saved workbook view
$ws.Range("D3").Select()
